# =====================================================================
# Apply "Added bruteforce algorithm and more test data" edit
#
# Original sheet had a single "Nearest Addition on Route A" block in
# A15:H38. We relocate that block down to A41:H64 (unchanged), and then
# fill in three new blocks of test data:
#   - Nearest Addition on Route B  -> J41:Q64
#   - Nearest Addition on Route C  -> A15:H38
#   - Nearest Addition on Route D  -> J15:Q38
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: relocate the existing "Route A" block's formatting down to
# rows 41-64 before we overwrite A15:H38 with "Route C" content. We
# copy formats only (so borders / fills move, values will be re-applied
# explicitly afterwards).
# ---------------------------------------------------------------------
$ws.Range("A16:H38").Copy() | Out-Null
$ws.Range("A42").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Step 2: duplicate that same formatting into the two new columns-block
# areas (J41:Q64 for Route B, J15:Q38 for Route D) so the new tables
# look like the existing ones (borders, "Length:" box, etc.).
# ---------------------------------------------------------------------
$ws.Range("A16:H38").Copy() | Out-Null
$ws.Range("J41").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A16:H38").Copy() | Out-Null
$ws.Range("J15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Step 3: re-write the moved "Route A" block's own values/formulas
# (content identical to before, just shifted down by 26 rows).
# ---------------------------------------------------------------------
$ws.Range("A41").Value = "Nearest Addition on Route A"
$ws.Range("A42").Value = "Route:"
$ws.Range("B42").Value = 0
$ws.Range("G42").Value = "Length:"
$ws.Range("H42").Value = 0
$ws.Range("A43").Value = "Route:"
$ws.Range("B43").Value = "0 5"
$ws.Range("D43").Value = "Add between 0 and 5"
$ws.Range("G43").Value = "Length:"
$ws.Range("H43").Formula = "=F7+F7"
$ws.Range("D44").Value = 2
$ws.Range("E44").Formula = "=F4+H7-F7"
$ws.Range("D45").Value = 4
$ws.Range("E45").Formula = "=F6+J7-F7"
$ws.Range("D46").Value = 9
$ws.Range("E46").Formula = "=F11+K11-F7"
$ws.Range("A47").Value = "Route:"
$ws.Range("B47").Value = "0 2 5"
$ws.Range("D47").Value = "Add between 0 and 2"
$ws.Range("G47").Value = "Length:"
$ws.Range("H47").Formula = "=F4+H7+F7"
$ws.Range("D48").Value = 4
$ws.Range("E48").Formula = "=F6+H6-F4"
$ws.Range("D49").Value = 9
$ws.Range("E49").Formula = "=F11+H11-F4"
$ws.Range("D50").Value = "Add between 2 and 5"
$ws.Range("D51").Value = 4
$ws.Range("E51").Formula = "=H6+J7-H7"
$ws.Range("D52").Value = 9
$ws.Range("E52").Formula = "=H11+K11-H7"
$ws.Range("D53").Value = "Add between 0 and 5"
$ws.Range("D54").Value = 4
$ws.Range("E54").Formula = "=F6+J7-F7"
$ws.Range("D55").Value = 9
$ws.Range("E55").Formula = "=F11+K11-F7"
$ws.Range("A56").Value = "Route:"
$ws.Range("B56").Value = "0 2 4 5"
$ws.Range("D56").Value = "Add between 0 and 2"
$ws.Range("G56").Value = "Length:"
$ws.Range("H56").Formula = "=F4+H6+J7+F7"
$ws.Range("D57").Value = 9
$ws.Range("E57").Formula = "=F11+H11-F4"
$ws.Range("D58").Value = "Add between 2 and 4"
$ws.Range("D59").Value = 9
$ws.Range("E59").Formula = "=H11+J11-H6"
$ws.Range("D60").Value = "Add between 4 and 5"
$ws.Range("D61").Value = 9
$ws.Range("E61").Formula = "=J11+K11-J7"
$ws.Range("D62").Value = "Add between 0 and 5"
$ws.Range("D63").Value = 9
$ws.Range("E63").Formula = "=F11+K11-F7"
$ws.Range("A64").Value = "Route:"
$ws.Range("B64").Value = "0 2 4 9 5"
$ws.Range("G64").Value = "Length:"
$ws.Range("H64").Formula = "=F4+H6+J11+K11+F7"

# ---------------------------------------------------------------------
# Step 4: new "Route B" data (J41:Q64).
# ---------------------------------------------------------------------
$ws.Range("J41").Value = "Nearest Addition on Route B"
$ws.Range("J42").Value = "Route:"
$ws.Range("K42").Value = 0
$ws.Range("J43").Value = "Route:"
$ws.Range("K43").Value = "0 6"
$ws.Range("M43").Value = "Add between 0 and 6"
$ws.Range("M44").Value = 3
$ws.Range("N44").Formula = "=F5+I8-F8"
$ws.Range("M45").Value = 7
$ws.Range("N45").Formula = "=F9+L9-F8"
$ws.Range("M46").Value = 10
$ws.Range("N46").Formula = "=F12+L12-F8"
$ws.Range("J47").Value = "Route:"
$ws.Range("K47").Value = "0 7 6"
$ws.Range("M47").Value = "Add between 0 and 7"
$ws.Range("M48").Value = 3
$ws.Range("N48").Formula = "=F5+I9-F9"
$ws.Range("M49").Value = 10
$ws.Range("N49").Formula = "=F12+M12-F9"
$ws.Range("M50").Value = "Add between 7 and 6"
$ws.Range("M51").Value = 3
$ws.Range("N51").Formula = "=I9+I8-L9"
$ws.Range("M52").Value = 10
$ws.Range("N52").Formula = "=M12+L12-L9"
$ws.Range("M53").Value = "Add between 0 and 6"
$ws.Range("M54").Value = 3
$ws.Range("N54").Formula = "=F5+I8-F8"
$ws.Range("M55").Value = 10
$ws.Range("N55").Formula = "=F12+L12-F8"
$ws.Range("J56").Value = "Route:"
$ws.Range("K56").Value = "0 3 7 6"
$ws.Range("M56").Value = "Add between 0 and 3"
$ws.Range("M57").Value = 10
$ws.Range("N57").Formula = "=F12+I12-F5"
$ws.Range("M58").Value = "Add between 3 and 7"
$ws.Range("M59").Value = 10
$ws.Range("N59").Formula = "=I12+M12-I9"
$ws.Range("M60").Value = "Add between 7 and 6"
$ws.Range("M61").Value = 10
$ws.Range("N61").Formula = "=M12+L12-L9"
$ws.Range("M62").Value = "Add between 0 and 6"
$ws.Range("M63").Value = 10
$ws.Range("N63").Formula = "=F12+L12-F8"
$ws.Range("J64").Value = "Route:"
$ws.Range("K64").Value = "0 3 7 10 6"
$ws.Range("P64").Value = "Length:"
$ws.Range("Q64").Formula = "=F5+I9+M12+L12+F8"

# ---------------------------------------------------------------------
# Step 5: new "Route C" data (A15:H38, replacing the old "Route A"
# content that used to live here).
# ---------------------------------------------------------------------
$ws.Range("A15").Value = "Nearest Addition on Route C"
$ws.Range("A16").Value = "Route:"
$ws.Range("B16").Value = 0
$ws.Range("A17").Value = "Route:"
$ws.Range("B17").Value = "0 1"
$ws.Range("D17").Value = "Add between 0 and 1"
$ws.Range("D18").Value = 2
$ws.Range("E18").Formula = "=F4+G4-F3"
$ws.Range("D19").Value = 7
$ws.Range("E19").Formula = "=F9+G9-F3"
$ws.Range("D20").Value = 10
$ws.Range("E20").Formula = "=F12+G12-F3"
$ws.Range("A21").Value = "Route:"
$ws.Range("B21").Value = "0 2 1 "
$ws.Range("D21").Value = "Add between 0 and 2"
$ws.Range("D22").Value = 7
$ws.Range("E22").Formula = "=F9+H9-F4"
$ws.Range("D23").Value = 10
$ws.Range("E23").Formula = "=F12+H12-F4"
$ws.Range("D24").Value = "Add between 2 and 1"
$ws.Range("D25").Value = 7
$ws.Range("E25").Formula = "=H9+G9-G4"
$ws.Range("D26").Value = 10
$ws.Range("E26").Formula = "=H12+G12-G4"
$ws.Range("D27").Value = "Add between 0 and 1"
$ws.Range("D28").Value = 7
$ws.Range("E28").Formula = "=F9+G9-F3"
$ws.Range("D29").Value = 10
$ws.Range("E29").Formula = "=F12+G12-F3"
$ws.Range("A30").Value = "Route:"
$ws.Range("B30").Value = "0 7 2 1"
$ws.Range("D30").Value = "Add between 0 and 7"
$ws.Range("D31").Value = 10
$ws.Range("E31").Formula = "=F12+M12-F9"
$ws.Range("D32").Value = "Add between 7 and 2"
$ws.Range("D33").Value = 10
$ws.Range("E33").Formula = "=M12+H12-H9"
$ws.Range("D34").Value = "Add between 2 and 1"
$ws.Range("D35").Value = 10
$ws.Range("E35").Formula = "=H12+G12-G4"
$ws.Range("D36").Value = "Add between 0 and 1"
$ws.Range("D37").Value = 10
$ws.Range("E37").Formula = "=F12+G12-F3"
$ws.Range("A38").Value = "Route:"
$ws.Range("B38").Value = "0 10 7 2 1"
$ws.Range("G38").Value = "Length:"
$ws.Range("H38").Formula = "=F12+M12+H9+G4+F3"

# ---------------------------------------------------------------------
# Step 6: new "Route D" data (J15:Q38).
# ---------------------------------------------------------------------
$ws.Range("J15").Value = "Nearest Addition on Route D"
$ws.Range("J16").Value = "Route:"
$ws.Range("K16").Value = 0
$ws.Range("J17").Value = "Route:"
$ws.Range("K17").Value = "0 5"
$ws.Range("M17").Value = "Add between 0 and 5"
$ws.Range("M18").Value = 2
$ws.Range("N18").Formula = "=F4+H7-F7"
$ws.Range("M19").Value = 4
$ws.Range("N19").Formula = "=F6+J7-F7"
$ws.Range("M20").Value = 8
$ws.Range("N20").Formula = "=F10+K10-F7"
$ws.Range("J21").Value = "Route:"
$ws.Range("K21").Value = "0 8 5"
$ws.Range("M21").Value = "Add between 0 and 8"
$ws.Range("M22").Value = 2
$ws.Range("N22").Formula = "=F4+H10-F10"
$ws.Range("M23").Value = 4
$ws.Range("N23").Formula = "=F6+J10-F10"
$ws.Range("M24").Value = "Add between 8 and 5"
$ws.Range("M25").Value = 2
$ws.Range("N25").Formula = "=H10+H7-K10"
$ws.Range("M26").Value = 4
$ws.Range("N26").Formula = "=J10+J7-K10"
$ws.Range("M27").Value = "Add between 0 and 5"
$ws.Range("M28").Value = 2
$ws.Range("N28").Formula = "=F4+H7-F7"
$ws.Range("M29").Value = 4
$ws.Range("N29").Formula = "=F6+J7-F7"
$ws.Range("J30").Value = "Route:"
$ws.Range("K30").Value = "0 4 8 5"
$ws.Range("M30").Value = "Add between 0 and 4"
$ws.Range("M31").Value = 2
$ws.Range("N31").Formula = "=F4+H6-F6"
$ws.Range("M32").Value = "Add between 4 and 8"
$ws.Range("M33").Value = 2
$ws.Range("N33").Formula = "=H6+H10-J10"
$ws.Range("M34").Value = "Add between 8 and 5"
$ws.Range("M35").Value = 2
$ws.Range("N35").Formula = "=H10+H7-K10"
$ws.Range("M36").Value = "Add between 0 and 5"
$ws.Range("M37").Value = 2
$ws.Range("N37").Formula = "=F4+H7-F7"
$ws.Range("J38").Value = "Route:"
$ws.Range("K38").Value = "0 2 4 8 5"
$ws.Range("P38").Value = "Length:"
$ws.Range("Q38").Formula = "=F4+H6+J10+K10+F7"

# ---------------------------------------------------------------------
# Step 7: highlight the "best addition" length cells (same green fill
# used by the existing table) in every block.
# ---------------------------------------------------------------------
$ws.Range("E18").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("N20").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("N23").PasteSpecial(-4122) | Out-Null
$ws.Range("E31").PasteSpecial(-4122) | Out-Null
$ws.Range("N31").PasteSpecial(-4122) | Out-Null
$ws.Range("E44").PasteSpecial(-4122) | Out-Null
$ws.Range("N45").PasteSpecial(-4122) | Out-Null
$ws.Range("N48").PasteSpecial(-4122) | Out-Null
$ws.Range("E51").PasteSpecial(-4122) | Out-Null
$ws.Range("E61").PasteSpecial(-4122) | Out-Null
$ws.Range("N61").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Step 8: view bookkeeping - scroll a bit right and move the selection,
# matching the author's final cursor position.
# ---------------------------------------------------------------------
$ws.Range("U1").Select() | Out-Null
try { $excel.ActiveWindow.ScrollColumn = 2 } catch { }
